$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 770776.1
$ws.Range("I2").Value = 946
$ws.Range("J2").Value = 1001725.2
$ws.Range("K2").Value = 946
$ws.Range("L2").Value = 1001725.2
$ws.Range("M2").Value = -833
$ws.Range("N2").Value = -1001951.2
$ws.Range("H53").Value = 445.8
$ws.Range("I53").Value = 557.5
$ws.Range("K53").Value = 557.5
$ws.Range("M53").Value = 79.5
$ws.Range("H69").Value = 23250
$ws.Range("J69").Value = 23250
$ws.Range("L69").Value = 69750
$ws.Range("N69").Value = -71498
$ws.Range("H72").Value = 23250
$ws.Range("J72").Value = 23250
$ws.Range("L72").Value = 209250
$ws.Range("N72").Value = -217986
$ws.Range("H80").Value = 2319196.5
$ws.Range("I80").Value = 1634551
$ws.Range("K80").Value = 4903653
$ws.Range("M80").Value = -4902655
$ws.Range("H83").Value = 2319196.5
$ws.Range("I83").Value = 1634551
$ws.Range("K83").Value = 14710959
$ws.Range("M83").Value = -14705967
$ws.Range("H92").Value = 1943.579
$ws.Range("I92").Value = 1045.2727
$ws.Range("K92").Value = 1045.2727
$ws.Range("M92").Value = 202.7273
$ws.Range("H98").Value = 5001589.5
$ws.Range("I98").Value = 5953779
$ws.Range("K98").Value = 5953779
$ws.Range("M98").Value = -5952281
$ws.Range("H100").Value = 10683.546
$ws.Range("J100").Value = 13314.875
$ws.Range("L100").Value = 13314.875
$ws.Range("N100").Value = -14396.875
$ws.Range("H107").Value = 866.72
$ws.Range("J107").Value = 1433.8889
$ws.Range("L107").Value = 1433.8889
$ws.Range("N107").Value = -5273.8889
$ws.Range("H116").Value = 7389.5
$ws.Range("I116").Value = 6649.5
$ws.Range("J116").Value = 7463.5
$ws.Range("K116").Value = 6649.5
$ws.Range("L116").Value = 7463.5
$ws.Range("M116").Value = -3207.5
$ws.Range("N116").Value = -14347.5
$ws.Range("H122").Value = 5001589.5
$ws.Range("I122").Value = 5953779
$ws.Range("K122").Value = 17861337
$ws.Range("M122").Value = -17858887
$ws.Range("H132").Value = 6288.115
$ws.Range("I132").Value = 3639.8
$ws.Range("K132").Value = 10919.4
$ws.Range("M132").Value = -8389.400000000001
$ws.Range("H134").Value = 39331.668
$ws.Range("J134").Value = 39331.668
$ws.Range("L134").Value = 39331.668
$ws.Range("N134").Value = -49471.668
$ws.Range("H137").Value = 2312.2307
$ws.Range("I137").Value = 2348.889
$ws.Range("J137").Value = 2229.75
$ws.Range("K137").Value = 7046.667
$ws.Range("L137").Value = 6689.25
$ws.Range("M137").Value = -4496.667
$ws.Range("N137").Value = -11789.25
$ws.Range("H138").Value = 7780.864
$ws.Range("I138").Value = 3896.7693
$ws.Range("J138").Value = 13391.223
$ws.Range("K138").Value = 11690.3079
$ws.Range("L138").Value = 40173.669
$ws.Range("M138").Value = -6550.3079
$ws.Range("N138").Value = -50453.669
$ws.Range("H141").Value = 15172556
$ws.Range("I141").Value = 21745718
$ws.Range("J141").Value = 54283.5
$ws.Range("K141").Value = 65237154
$ws.Range("L141").Value = 162850.5
$ws.Range("M141").Value = -65231974
$ws.Range("N141").Value = -173210.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4337.4375
$ws.Range("I2").Value = 4476
$ws.Range("J2").Value = 4159.2856
$ws.Range("K2").Value = 4476
$ws.Range("L2").Value = 4159.2856
$ws.Range("M2").Value = -4363
$ws.Range("N2").Value = -4385.2856
$ws.Range("H32").Value = 3174.0715
$ws.Range("I32").Value = 2533.2834
$ws.Range("J32").Value = 7018.8
$ws.Range("K32").Value = 2533.2834
$ws.Range("L32").Value = 7018.8
$ws.Range("M32").Value = -2246.2834
$ws.Range("N32").Value = -7592.8
$ws.Range("H39").Value = 16507.5
$ws.Range("I39").Value = 16507.5
$ws.Range("K39").Value = 16507.5
$ws.Range("M39").Value = -15987.5
$ws.Range("H45").Value = 3284.2856
$ws.Range("I45").Value = 2502.5
$ws.Range("J45").Value = 4326.6665
$ws.Range("K45").Value = 2502.5
$ws.Range("L45").Value = 4326.6665
$ws.Range("M45").Value = -2125.5
$ws.Range("N45").Value = -5080.6665
$ws.Range("H61").Value = 25556918
$ws.Range("I61").Value = 52501228
$ws.Range("K61").Value = 52501228
$ws.Range("M61").Value = -52501016
$ws.Range("H74").Value = 2278.1155
$ws.Range("I74").Value = 1331.1578
$ws.Range("K74").Value = 1331.1578
$ws.Range("M74").Value = -457.1578
$ws.Range("H77").Value = 2278.1155
$ws.Range("I77").Value = 1331.1578
$ws.Range("K77").Value = 6655.789
$ws.Range("M77").Value = -2287.789
$ws.Range("H110").Value = 1779.5
$ws.Range("I110").Value = 706.3333
$ws.Range("K110").Value = 706.3333
$ws.Range("M110").Value = 1338.6667
$ws.Range("H116").Value = 4337.4375
$ws.Range("I116").Value = 4476
$ws.Range("J116").Value = 4159.2856
$ws.Range("K116").Value = 4476
$ws.Range("L116").Value = 4159.2856
$ws.Range("M116").Value = -2182
$ws.Range("N116").Value = -8747.285599999999
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 3035946.8
$ws.Range("I132").Value = 5123.5
$ws.Range("J132").Value = 11118142
$ws.Range("K132").Value = 15370.5
$ws.Range("L132").Value = 33354426
$ws.Range("M132").Value = -12840.5
$ws.Range("N132").Value = -33359486
$ws.Range("H136").Value = 25556918
$ws.Range("I136").Value = 52501228
$ws.Range("K136").Value = 157503684
$ws.Range("M136").Value = -157501134
$ws.Range("H141").Value = 25000
$ws.Range("I141").Value = 25000
$ws.Range("K141").Value = 25000
$ws.Range("M141").Value = -19820

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4337.4375
$ws.Range("I3").Value = 4476
$ws.Range("J3").Value = 4159.2856
$ws.Range("K3").Value = 4476
$ws.Range("L3").Value = 4159.2856
$ws.Range("M3").Value = -4362
$ws.Range("N3").Value = -4387.2856
$ws.Range("H20").Value = 7324.467
$ws.Range("J20").Value = 1621.4
$ws.Range("L20").Value = 1621.4
$ws.Range("N20").Value = -2115.4
$ws.Range("H82").Value = 27993.5
$ws.Range("J82").Value = 99999
$ws.Range("L82").Value = 99999
$ws.Range("N82").Value = -100765
$ws.Range("H85").Value = 27993.5
$ws.Range("J85").Value = 99999
$ws.Range("L85").Value = 99999
$ws.Range("N85").Value = -102651
$ws.Range("H86").Value = 43892.293
$ws.Range("I86").Value = 68729
$ws.Range("K86").Value = 68729
$ws.Range("M86").Value = -67606
$ws.Range("H89").Value = 43892.293
$ws.Range("I89").Value = 68729
$ws.Range("K89").Value = 343645
$ws.Range("M89").Value = -338029
$ws.Range("H94").Value = 1573.1923
$ws.Range("I94").Value = 1551.9048
$ws.Range("K94").Value = 1551.9048
$ws.Range("M94").Value = -1100.9048
$ws.Range("H99").Value = 3988
$ws.Range("I99").Value = 3988
$ws.Range("K99").Value = 3988
$ws.Range("M99").Value = -2490
$ws.Range("H105").Value = 807110.25
$ws.Range("I105").Value = 1608158.2
$ws.Range("J105").Value = 6062.25
$ws.Range("K105").Value = 1608158.2
$ws.Range("L105").Value = 6062.25
$ws.Range("M105").Value = -1606411.2
$ws.Range("N105").Value = -9556.25
$ws.Range("H107").Value = 2295.923
$ws.Range("I107").Value = 2168.182
$ws.Range("K107").Value = 2168.182
$ws.Range("M107").Value = -248.1819999999998
$ws.Range("H134").Value = 4764660.5
$ws.Range("I134").Value = 2666.0667
$ws.Range("J134").Value = 16669647
$ws.Range("K134").Value = 7998.2001
$ws.Range("L134").Value = 50008941
$ws.Range("M134").Value = -5463.2001
$ws.Range("N134").Value = -50014011

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 31253306
$ws.Range("I31").Value = 55558844
$ws.Range("J31").Value = 3325.4285
$ws.Range("K31").Value = 55558844
$ws.Range("L31").Value = 3325.4285
$ws.Range("M31").Value = -55558549
$ws.Range("N31").Value = -3915.4285
$ws.Range("H34").Value = 31253306
$ws.Range("I34").Value = 55558844
$ws.Range("J34").Value = 3325.4285
$ws.Range("K34").Value = 55558844
$ws.Range("L34").Value = 3325.4285
$ws.Range("M34").Value = -55558642
$ws.Range("N34").Value = -3729.4285
$ws.Range("H35").Value = 4397.5713
$ws.Range("I35").Value = 4642.5
$ws.Range("K35").Value = 4642.5
$ws.Range("M35").Value = -4348.5
$ws.Range("H58").Value = 3012.7856
$ws.Range("J58").Value = 2999.75
$ws.Range("L58").Value = 2999.75
$ws.Range("N58").Value = -3405.75
$ws.Range("H94").Value = 941.5294
$ws.Range("I94").Value = 845.375
$ws.Range("K94").Value = 845.375
$ws.Range("M94").Value = -394.375
$ws.Range("H97").Value = 99999
$ws.Range("J97").Value = 99999
$ws.Range("L97").Value = 99999
$ws.Range("N97").Value = -101981
$ws.Range("H99").Value = 21714.857
$ws.Range("I99").Value = 25000.666
$ws.Range("K99").Value = 25000.666
$ws.Range("M99").Value = -23502.666
$ws.Range("H105").Value = 1892.8182
$ws.Range("I105").Value = 1349.6471
$ws.Range("J105").Value = 3739.6
$ws.Range("K105").Value = 1349.6471
$ws.Range("L105").Value = 3739.6
$ws.Range("M105").Value = 397.3529000000001
$ws.Range("N105").Value = -7233.6
$ws.Range("H107").Value = 1040.1034
$ws.Range("I107").Value = 745.25
$ws.Range("J107").Value = 2455.4
$ws.Range("K107").Value = 745.25
$ws.Range("L107").Value = 2455.4
$ws.Range("M107").Value = 1174.75
$ws.Range("N107").Value = -6295.4
$ws.Range("H126").Value = 21714.857
$ws.Range("I126").Value = 25000.666
$ws.Range("K126").Value = 75001.99800000001
$ws.Range("M126").Value = -72531.99800000001
$ws.Range("H132").Value = 3505.625
$ws.Range("I132").Value = 3509
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 10527
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -7997
$ws.Range("N132").Value = -15560
$ws.Range("H134").Value = 3889.125
$ws.Range("I134").Value = 3683.1667
$ws.Range("J134").Value = 4507
$ws.Range("K134").Value = 11049.5001
$ws.Range("L134").Value = 13521
$ws.Range("M134").Value = -8514.500100000001
$ws.Range("N134").Value = -18591
$ws.Range("H136").Value = 3012.7856
$ws.Range("J136").Value = 2999.75
$ws.Range("L136").Value = 8999.25
$ws.Range("N136").Value = -14099.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 79.28570999999999
$ws.Range("J2").Value = 88.8
$ws.Range("L2").Value = 532.8
$ws.Range("N2").Value = -758.8
$ws.Range("H11").Value = 5208.909
$ws.Range("I11").Value = 1479.6666
$ws.Range("K11").Value = 4438.9998
$ws.Range("M11").Value = -4298.9998
$ws.Range("H12").Value = 1314.0714
$ws.Range("I12").Value = 241.16667
$ws.Range("J12").Value = 2118.75
$ws.Range("K12").Value = 723.50001
$ws.Range("L12").Value = 6356.25
$ws.Range("M12").Value = -550.50001
$ws.Range("N12").Value = -6702.25
$ws.Range("H26").Value = 6826
$ws.Range("I26").Value = 199.66667
$ws.Range("J26").Value = 16765.5
$ws.Range("K26").Value = 599.00001
$ws.Range("L26").Value = 50296.5
$ws.Range("M26").Value = -311.00001
$ws.Range("N26").Value = -50872.5
$ws.Range("H38").Value = 305.33334
$ws.Range("I38").Value = 7.75
$ws.Range("J38").Value = 900.5
$ws.Range("K38").Value = 23.25
$ws.Range("L38").Value = 2701.5
$ws.Range("M38").Value = 323.75
$ws.Range("N38").Value = -3395.5
$ws.Range("H44").Value = 6594.2856
$ws.Range("I44").Value = 226
$ws.Range("J44").Value = 11370.5
$ws.Range("K44").Value = 678
$ws.Range("L44").Value = 34111.5
$ws.Range("M44").Value = -280
$ws.Range("N44").Value = -34907.5
$ws.Range("H127").Value = 3799.6667
$ws.Range("J127").Value = 3799.6667
$ws.Range("L127").Value = 11399.0001
$ws.Range("N127").Value = -21319.0001
$ws.Range("H131").Value = 3832.6086
$ws.Range("I131").Value = 2136
$ws.Range("J131").Value = 6038.2
$ws.Range("K131").Value = 6408
$ws.Range("L131").Value = 18114.6
$ws.Range("M131").Value = -1368
$ws.Range("N131").Value = -28194.6
$ws.Range("H139").Value = 8173.8184
$ws.Range("I139").Value = 6131.1113
$ws.Range("K139").Value = 18393.3339
$ws.Range("M139").Value = -13253.3339
$ws.Range("H140").Value = 4089.4167
$ws.Range("I140").Value = 1430.909
$ws.Range("K140").Value = 4292.727000000001
$ws.Range("M140").Value = 887.2729999999992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2548.818
$ws.Range("J80").Value = 4200.5
$ws.Range("L80").Value = 4200.5
$ws.Range("N80").Value = -6196.5
$ws.Range("H83").Value = 2548.818
$ws.Range("J83").Value = 4200.5
$ws.Range("L83").Value = 21002.5
$ws.Range("N83").Value = -30986.5
$ws.Range("H97").Value = 786.875
$ws.Range("I97").Value = 813.5714
$ws.Range("K97").Value = 813.5714
$ws.Range("M97").Value = -317.5714
$ws.Range("H102").Value = 4237.5454
$ws.Range("I102").Value = 4161.4
$ws.Range("K102").Value = 4161.4
$ws.Range("M102").Value = -2539.4
$ws.Range("H113").Value = 1160131
$ws.Range("I113").Value = 2821.5557
$ws.Range("J113").Value = 2648100.2
$ws.Range("K113").Value = 2821.5557
$ws.Range("L113").Value = 2648100.2
$ws.Range("M113").Value = -651.5556999999999
$ws.Range("N113").Value = -2652440.2
$ws.Range("H122").Value = 1341.6666
$ws.Range("I122").Value = 1341.6666
$ws.Range("K122").Value = 4024.9998
$ws.Range("M122").Value = -1574.9998
$ws.Range("H132").Value = 7146603.5
$ws.Range("I132").Value = 4026.5
$ws.Range("J132").Value = 25003046
$ws.Range("K132").Value = 12079.5
$ws.Range("L132").Value = 75009138
$ws.Range("M132").Value = -9549.5
$ws.Range("N132").Value = -75014198

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6423.222
$ws.Range("I7").Value = 6663.625
$ws.Range("J7").Value = 4500
$ws.Range("K7").Value = 6663.625
$ws.Range("L7").Value = 4500
$ws.Range("M7").Value = -6551.625
$ws.Range("N7").Value = -4724
$ws.Range("H40").Value = 2913.0476
$ws.Range("I40").Value = 3071.889
$ws.Range("K40").Value = 3071.889
$ws.Range("M40").Value = -2935.889
$ws.Range("H55").Value = 926.05884
$ws.Range("I55").Value = 468.6875
$ws.Range("J55").Value = 1332.6111
$ws.Range("K55").Value = 468.6875
$ws.Range("L55").Value = 1332.6111
$ws.Range("M55").Value = -295.6875
$ws.Range("N55").Value = -1678.6111
$ws.Range("H61").Value = 83338980
$ws.Range("I61").Value = 166667620
$ws.Range("K61").Value = 166667620
$ws.Range("M61").Value = -166667418
$ws.Range("H82").Value = 3768
$ws.Range("I82").Value = 2563.6
$ws.Range("K82").Value = 2563.6
$ws.Range("M82").Value = -2202.6
$ws.Range("H85").Value = 3768
$ws.Range("I85").Value = 2563.6
$ws.Range("K85").Value = 2563.6
$ws.Range("M85").Value = -1315.6
$ws.Range("H93").Value = 6179781
$ws.Range("I93").Value = 4412.5
$ws.Range("K93").Value = 4412.5
$ws.Range("M93").Value = -3164.5
$ws.Range("H113").Value = 83338980
$ws.Range("I113").Value = 166667620
$ws.Range("K113").Value = 166667620
$ws.Range("M113").Value = -166665450
$ws.Range("H122").Value = 3533.7173
$ws.Range("I122").Value = 3256.0715
$ws.Range("J122").Value = 6449
$ws.Range("K122").Value = 9768.2145
$ws.Range("L122").Value = 19347
$ws.Range("M122").Value = -7318.2145
$ws.Range("N122").Value = -24247
$ws.Range("H125").Value = 91249.75
$ws.Range("J125").Value = 91249.75
$ws.Range("L125").Value = 91249.75
$ws.Range("N125").Value = -101089.75
$ws.Range("H126").Value = 6423.222
$ws.Range("I126").Value = 6663.625
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 19990.875
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -17520.875
$ws.Range("N126").Value = -18440
$ws.Range("H132").Value = 3538.3044
$ws.Range("I132").Value = 2292.8125
$ws.Range("J132").Value = 6385.143
$ws.Range("K132").Value = 6878.4375
$ws.Range("L132").Value = 19155.429
$ws.Range("M132").Value = -4348.4375
$ws.Range("N132").Value = -24215.429
$ws.Range("H136").Value = 4017.5625
$ws.Range("I136").Value = 3852.0667
$ws.Range("K136").Value = 11556.2001
$ws.Range("M136").Value = -9006.2001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H82").Value = 71929.39999999999
$ws.Range("J82").Value = 71929.39999999999
$ws.Range("L82").Value = 71929.39999999999
$ws.Range("N82").Value = -72695.39999999999
$ws.Range("H85").Value = 71929.39999999999
$ws.Range("J85").Value = 71929.39999999999
$ws.Range("L85").Value = 71929.39999999999
$ws.Range("N85").Value = -74581.39999999999
$ws.Range("H107").Value = 4229.968
$ws.Range("J107").Value = 4642.5415
$ws.Range("L107").Value = 13927.6245
$ws.Range("N107").Value = -17767.6245
$ws.Range("H110").Value = 120000
$ws.Range("J110").Value = 120000
$ws.Range("L110").Value = 120000
$ws.Range("N110").Value = -128180
$ws.Range("H122").Value = 2438.6667
$ws.Range("I122").Value = 2431.0625
$ws.Range("K122").Value = 7293.1875
$ws.Range("M122").Value = -4843.1875
$ws.Range("H132").Value = 189477.27
$ws.Range("I132").Value = 4151.0713
$ws.Range("J132").Value = 838118.9399999999
$ws.Range("K132").Value = 12453.2139
$ws.Range("L132").Value = 2514356.82
$ws.Range("M132").Value = -9923.213899999999
$ws.Range("N132").Value = -2519416.82
$ws.Range("H136").Value = 452021.6
$ws.Range("I136").Value = 17976.047
$ws.Range("K136").Value = 53928.141
$ws.Range("M136").Value = -51378.141
